# Applies the "Trade #95 closed" update to live_trading_results.xlsx
#  - Summary: Total Trades 122->123, Win Rate % 47.54->47.15
#  - Strategy Status: MarketMaking row Trades 42->43, Win Rate % 47.62->46.51
#  - All Trades: closes existing MarketMaking OPEN trade (row 125) and
#    appends two new OPEN trades (momentum, HighProbConvergence)
#  - momentum sheet: appends the new OPEN momentum trade
#  - HighProbConvergence sheet: appends the new OPEN HighProbConvergence trade
#  - MarketMaking sheet: closes the existing OPEN trade (row 45)

$wb = $excel.ActiveWorkbook

# Helper: write a text value into a cell without Excel's automatic
# date/number inference turning things like "2026-02-18" into a date serial.
function Set-TextCell($ws, $addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 123    # Total Trades
$summary.Range("B9").Value = 47.15  # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking is row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D6").Value = 43     # Trades
$status.Range("G6").Value = 46.51  # Win Rate %

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close the existing OPEN MarketMaking trade on row 125 (Trade # 124)
$allTrades.Range("G125").Value = 0.96
$allTrades.Range("H125").Value = "CLOSED"
$allTrades.Range("K125").Value = 99.53
$allTrades.Range("L125").Value = "early_exit"
$allTrades.Range("M125").Value = 0.1

# New row 154: Trade # 153, momentum strategy, still OPEN
$allTrades.Range("A154").Value = 153
Set-TextCell $allTrades "B154" "2026-02-18"
Set-TextCell $allTrades "C154" "00:33:34"
Set-TextCell $allTrades "D154" "momentum"
Set-TextCell $allTrades "E154" "UP"
$allTrades.Range("F154").Value = 0.96
Set-TextCell $allTrades "H154" "OPEN"
$allTrades.Range("I154").Value = 0
$allTrades.Range("J154").Value = 0
$allTrades.Range("K154").Value = 99.23374292899115
$allTrades.Range("M154").Value = 0
$allTrades.Range("N154").Value = 0
$allTrades.Range("O154").Value = 0
$allTrades.Range("P154").Value = 0.9
Set-TextCell $allTrades "Q154" "Upward momentum: 1.980% over 10 samples"

# New row 155: Trade # 154, HighProbConvergence strategy, still OPEN
$allTrades.Range("A155").Value = 154
Set-TextCell $allTrades "B155" "2026-02-18"
Set-TextCell $allTrades "C155" "00:33:35"
Set-TextCell $allTrades "D155" "HighProbConvergence"
Set-TextCell $allTrades "E155" "DOWN"
$allTrades.Range("F155").Value = 0.04
Set-TextCell $allTrades "H155" "OPEN"
$allTrades.Range("I155").Value = 0
$allTrades.Range("J155").Value = 0
$allTrades.Range("K155").Value = 100.4130057263667
$allTrades.Range("M155").Value = 0
$allTrades.Range("N155").Value = 0
$allTrades.Range("O155").Value = 0
$allTrades.Range("P155").Value = 0.95
Set-TextCell $allTrades "Q155" "Mean reversion DOWN: price 1.88% above mean (z=4.36)"

# ---------------------------------------------------------------------
# momentum sheet: append new row 38 (Trade # 153)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("A38").Value = 153
Set-TextCell $momentum "B38" "2026-02-18"
Set-TextCell $momentum "C38" "00:33:34"
Set-TextCell $momentum "D38" "momentum"
Set-TextCell $momentum "E38" "UP"
$momentum.Range("F38").Value = 0.96
Set-TextCell $momentum "H38" "OPEN"
$momentum.Range("I38").Value = 0
$momentum.Range("J38").Value = 0
$momentum.Range("K38").Value = 99.23374292899115
$momentum.Range("L38").Value = 0
$momentum.Range("M38").Value = 0
$momentum.Range("N38").Value = 0.9
Set-TextCell $momentum "O38" "Upward momentum: 1.980% over 10 samples"
$momentum.Range("Q38").Value = 0

# ---------------------------------------------------------------------
# HighProbConvergence sheet: append new row 20 (Trade # 154)
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Range("A20").Value = 154
Set-TextCell $hpc "B20" "2026-02-18"
Set-TextCell $hpc "C20" "00:33:35"
Set-TextCell $hpc "D20" "HighProbConvergence"
Set-TextCell $hpc "E20" "DOWN"
$hpc.Range("F20").Value = 0.04
Set-TextCell $hpc "H20" "OPEN"
$hpc.Range("I20").Value = 0
$hpc.Range("J20").Value = 0
$hpc.Range("K20").Value = 100.4130057263667
$hpc.Range("L20").Value = 0
$hpc.Range("M20").Value = 0
$hpc.Range("N20").Value = 0.95
Set-TextCell $hpc "O20" "Mean reversion DOWN: price 1.88% above mean (z=4.36)"
$hpc.Range("Q20").Value = 0

# ---------------------------------------------------------------------
# MarketMaking sheet: close the existing OPEN trade on row 45 (Trade # 124)
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G45").Value = 0.96
$marketMaking.Range("H45").Value = "CLOSED"
$marketMaking.Range("K45").Value = 99.53
$marketMaking.Range("P45").Value = "early_exit"
$marketMaking.Range("Q45").Value = 0.1
